$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 (shifts rows 13-21 down to 14-22, preserving their
#     heights/content) to make room for the "Docentes responsaveis" value row ---
$ws.Rows("13").Insert()

# The new row 13 has no pre-existing per-column style (B/C default to style 1
# on write), so copy the column B/C formatting from row 14 (which still holds
# the original style 2 / style 3 formatting after the shift) before writing
# values, then drop the stray A13 cell that Insert left behind.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()

$ws.Range("B13").Value = '3403572 - Ismael Maciel de Mancilha'
$ws.Range("C13").Value = '3403572 - Ismael Maciel de Mancilha'

# --- Objetivos: (row 10) now gets the real objectives text instead of the
#     misplaced professor name ---
$ws.Range("B10").Value = 'Levar aos estudantes conhecimentos básicos sobre a Ciência e Tecnologia dos Alimentos, abrangendo aspectos sobre as instituições envolvidas, bem como a habilitação dos profissionais nas respectivas áreas de concentração da C&T de Alimentos; segurança alimentar / controle de qualidade; legislação vigente; matérias primas, técnicas de processamento; composição dos alimentos, novos produtos, características do setor alimentício e sua relação com o meio ambiente. Neste contexto, enfoca-se a evolução do desenvolvimento da C&T dos Alimentos com vistas a atender a demanda por alimentos de qualidade.'
$ws.Range("C10").Value = 'Levar aos estudantes conhecimentos básicos sobre a Ciência e Tecnologia dos Alimentos, abrangendo aspectos sobre as instituições envolvidas, bem como a habilitação dos profissionais nas respectivas áreas de concentração da C&T de Alimentos; segurança alimentar / controle de qualidade; legislação vigente; matérias primas, técnicas de processamento; composição dos alimentos, novos produtos, características do setor alimentício e sua relação com o meio ambiente. Neste contexto, enfoca-se a evolução do desenvolvimento da C&T dos Alimentos com vistas a atender a demanda por alimentos de qualidade.'

# --- Programa resumido: (row 14 after shift) gets the real short-syllabus text
#     instead of "Semestral" ---
$ws.Range("B14").Value = 'Ciência e Tecnologia de Alimentos: conceitos e objetivos (consumo de alimentos, segurança alimentar, integridade e características nutritivas dos alimentos); Generalidades sobre o Setor Alimentício Legislação; Matérias-Primas; Técnicas de Processamento: Composição dos Alimentos; Novos produtos'
$ws.Range("C14").Value = 'Ciência e Tecnologia de Alimentos: conceitos e objetivos (consumo de alimentos, segurança alimentar, integridade e características nutritivas dos alimentos); Generalidades sobre o Setor Alimentício Legislação; Matérias-Primas; Técnicas de Processamento: Composição dos Alimentos; Novos produtos'

# --- Programa: (row 16 after shift) gets the real syllabus text instead of a
#     misplaced date ---
$ws.Range("B16").Value = 'Introdução: conceitos de C&T de alimentos; áreas de concentração e respectivas atribuições dos profissionais; instituições envolvidas com o desenvolvimento da C&T Alimentos; embalagens ativas/inteligentes Segurança Alimentar: fatores envolvidos na toxi-infecção alimentar; alimentos orgânicos e Boas Práticas de Fabricação (BPF); microbiologia de alimentos; higiene industrial; análise de perigo dos pontos críticos de controle - APPCC. Matérias-Primas: importância e características dos segmentos produtores de matérias primas de origem animal, vegetal, microbiana e aditivos/ingredientes. Processamento de Alimentos: objetivos e caracterização dos diferentes métodos de processamento dos alimentos abrangendo técnicas de conservação, transformação e melhoria da qualidade; usos e aplicações de aditivos/ingredientes em alimentos. Novos Produtos: fatores que devem ser considerados no lançamento de novos produtos alimentícios. Composição dos Alimentos: composição e significância dos diferentes compostos encontrados nos alimentos focando em suas propriedades e funções'
$ws.Range("C16").Value = 'Introdução: conceitos de C&T de alimentos; áreas de concentração e respectivas atribuições dos profissionais; instituições envolvidas com o desenvolvimento da C&T Alimentos; embalagens ativas/inteligentes Segurança Alimentar: fatores envolvidos na toxi-infecção alimentar; alimentos orgânicos e Boas Práticas de Fabricação (BPF); microbiologia de alimentos; higiene industrial; análise de perigo dos pontos críticos de controle - APPCC. Matérias-Primas: importância e características dos segmentos produtores de matérias primas de origem animal, vegetal, microbiana e aditivos/ingredientes. Processamento de Alimentos: objetivos e caracterização dos diferentes métodos de processamento dos alimentos abrangendo técnicas de conservação, transformação e melhoria da qualidade; usos e aplicações de aditivos/ingredientes em alimentos. Novos Produtos: fatores que devem ser considerados no lançamento de novos produtos alimentícios. Composição dos Alimentos: composição e significância dos diferentes compostos encontrados nos alimentos focando em suas propriedades e funções'

# --- Metodo: / Criterio: / Norma de recuperacao: / Bibliografia: (rows 19-22
#     after shift) each shift up to the correct value ---
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

$ws.Range("B20").Value = '2 provas (P1 + P2), sendo que a NF = (P1 + P2) / 2'
$ws.Range("C20").Value = '2 provas (P1 + P2), sendo que a NF = (P1 + P2) / 2'

$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'

$ws.Range("B22").Value = 'POTTER, N.N., HOTCHKISS, J.H., Food Science - 5ª Ed. Chapman & Hall, 1995.ARAÚJO, J.M.A., Química de Alimentos: Teoria e Prática. 3ª Ed. Editora UFV, 2004.FRAZIER, W.C., Microbiologia de los Alimentos. Editora Acribia, Zaragoza-Epanha, 1981.EVANGELISTA, J., Tecnologia de Alimentos, Livraria Atheneu, RJ. 1987.FENEMA, O.R.,Principles of Food Science: Part I: Food Chemistry. Marcel Dekker, Inc. NY-USA.1975.BENWART, G.J., Basic Food Microbiology. AVI Publishing Company Inc. USA,1970.'
$ws.Range("C22").Value = 'POTTER, N.N., HOTCHKISS, J.H., Food Science - 5ª Ed. Chapman & Hall, 1995.ARAÚJO, J.M.A., Química de Alimentos: Teoria e Prática. 3ª Ed. Editora UFV, 2004.FRAZIER, W.C., Microbiologia de los Alimentos. Editora Acribia, Zaragoza-Epanha, 1981.EVANGELISTA, J., Tecnologia de Alimentos, Livraria Atheneu, RJ. 1987.FENEMA, O.R.,Principles of Food Science: Part I: Food Chemistry. Marcel Dekker, Inc. NY-USA.1975.BENWART, G.J., Basic Food Microbiology. AVI Publishing Company Inc. USA,1970.'

# --- Column A's width/style range no longer spans column B (column B now has
#     its own width/style definition), so split the 1:2 range down to just
#     column A ---
$ws.Columns("A").ColumnWidth = 30.7109375
